$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Fila 52: "Estetica - botones imagenes" -> asignar responsable y 50% de avance ---
$ws.Cells.Item(52, 2).Value = "Agustina"
$ws.Cells.Item(52, 3).Value = 0.5
$ws.Cells.Item(52, 3).NumberFormat = "0%"

# --- Fila 53: "Estetica - etiquetas..." -> asignar responsable y 100% de avance ---
$ws.Cells.Item(53, 2).Value = "Agustina"
$ws.Cells.Item(53, 3).Value = 1
$ws.Cells.Item(53, 3).NumberFormat = "0%"

# --- Fila 56: "Reunion con Ivan y Josefina - consultas" -> 100% de avance ---
$ws.Cells.Item(56, 3).Value = 1
$ws.Cells.Item(56, 3).NumberFormat = "0%"

# --- Fila 58: "Agregar los comandos abajo..." -> pasa de "en proceso" a 100% ---
$ws.Cells.Item(58, 3).Value = 1
$ws.Cells.Item(58, 3).NumberFormat = "0%"

# --- Fila 61: "No asignar cliente en venta de factura B" -> 100% de avance ---
$ws.Cells.Item(61, 3).Value = 1
$ws.Cells.Item(61, 3).NumberFormat = "0%"

# --- Nuevas tareas agregadas al final del listado ---
$ws.Cells.Item(62, 1).Value = "Error en iva cuando consumidor final en ticket"
$ws.Cells.Item(62, 2).Value = "Lucas"

$ws.Cells.Item(63, 1).Value = "Error en listado de control - pagos de la fecha"
$ws.Cells.Item(63, 2).Value = "Agustina"

$ws.Cells.Item(64, 1).Value = "Migracion de datos"
$ws.Cells.Item(64, 2).Value = "Agustina"

$ws.Cells.Item(65, 1).Value = "Paginado de consultas"

$ws.Cells.Item(66, 1).Value = "Reporte listado de cliente"

$ws.Cells.Item(67, 1).Value = "Reporte para contador"
$ws.Cells.Item(67, 2).Value = "Agustina"

# --- Refleja la selección final del usuario (última celda tocada en la edición) ---
$ws.Range("B68").Select()
